$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells so they match
# (bold font, border, centered alignment - same style index as other headers)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for columns I (I0) and J (IF)
$values = @{
    2  = @(7, 7)
    3  = @(7, 8)
    4  = @(8, 8)
    5  = @(7, 7)
    6  = @(6, 7)
    7  = @(8, 8)
    8  = @(8, 8)
    9  = @(7, 7)
    10 = @(7, 7)
    11 = @(5, 5)
    12 = @(5, 6)
    13 = @(5, 5)
    14 = @(8, 8)
    15 = @(4, 4)
    16 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
